$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.4
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3.1
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.83
$ws.Range("X2").Value = 11
$ws.Range("AH2").Value = 15
$ws.Range("AK2").Value = 26
$ws.Range("AM2").Value = 4.33
$ws.Range("AN2").Value = 13
$ws.Range("AP2").Value = 41
